# Repair ESS1 and PV1 relay settings on the "relays" sheet.
# Relay# 2 (row 3) and Relay# 3 (row 4) had stale 13.8kV-feeder-style
# settings; fix them to the correct 480V ESS1/PV1 relay settings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("relays")

# Relay #2 (row 3) -- ESS1
$ws.Range("B3").Value = 480
$ws.Range("C3").Value = 10
$ws.Range("G3").Value = 3500
$ws.Range("K3").Value = 1
$ws.Range("R3").Value = 0.40093768693724008

# Relay #3 (row 4) -- PV1
$ws.Range("B4").Value = 480
$ws.Range("C4").Value = 10
$ws.Range("G4").Value = 3500
$ws.Range("K4").Value = 2
$ws.Range("R4").Value = 0.40093768693724008

# Leave the selection on K3, matching where the edits were made.
$ws.Range("K3").Select()
